$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 updates
$ws.Range("A1").Value = 123456
$ws.Range("B1").Value = "test"
$ws.Range("C1").Value = 130
$ws.Range("D1").Value = 1

# Row 2 (new)
$ws.Range("A2").Value = 123654
$ws.Range("B2").Value = "test 3"
$ws.Range("C2").Value = 260
$ws.Range("D2").Value = 1

# Column width adjustments
# (ColumnWidth -> stored OOXML width has a fixed +5/6 padding offset in this
# runtime, so subtract it to land on the exact target stored widths of 7 and 4)
$ws.Range("B1").ColumnWidth = 7 - 5/6
$ws.Range("C1").ColumnWidth = 4 - 5/6
